$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 0.9848484848484849
$ws.Range("D8").Value = 0.95
$ws.Range("F8").Value = 0.9285714285714286
$ws.Range("G8").Value = 0.9880952380952381
$ws.Range("H8").Value = 0.962962962962963
